$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2080
$ws.Range("I51").Value = 1800
$ws.Range("J51").Value = 2500
$ws.Range("K51").Value = 1800
$ws.Range("L51").Value = 2500
$ws.Range("M51").Value = -1316
$ws.Range("N51").Value = -3468
$ws.Range("H129").Value = 869.4123499999999
$ws.Range("I129").Value = 441.41666
$ws.Range("J129").Value = 929.83527
$ws.Range("K129").Value = 1324.24998
$ws.Range("L129").Value = 2789.50581
$ws.Range("M129").Value = 3675.75002
$ws.Range("N129").Value = -12789.50581
$ws.Range("H137").Value = 1353.303
$ws.Range("I137").Value = 1218.5217
$ws.Range("J137").Value = 1663.3
$ws.Range("K137").Value = 3655.5651
$ws.Range("L137").Value = 4989.9
$ws.Range("M137").Value = -1105.5651
$ws.Range("N137").Value = -10089.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7083.7896
$ws.Range("I61").Value = 7366.222
$ws.Range("K61").Value = 7366.222
$ws.Range("M61").Value = -7154.222
$ws.Range("H74").Value = 1402.625
$ws.Range("I74").Value = 1408.579
$ws.Range("J74").Value = 1380
$ws.Range("K74").Value = 1408.579
$ws.Range("L74").Value = 1380
$ws.Range("M74").Value = -534.579
$ws.Range("N74").Value = -3128
$ws.Range("H77").Value = 1402.625
$ws.Range("I77").Value = 1408.579
$ws.Range("J77").Value = 1380
$ws.Range("K77").Value = 7042.895
$ws.Range("L77").Value = 6900
$ws.Range("M77").Value = -2674.895
$ws.Range("N77").Value = -15636
$ws.Range("H132").Value = 3023.442
$ws.Range("I132").Value = 1687.4062
$ws.Range("J132").Value = 6910.091
$ws.Range("K132").Value = 5062.2186
$ws.Range("L132").Value = 20730.273
$ws.Range("M132").Value = -2532.2186
$ws.Range("N132").Value = -25790.273
$ws.Range("H136").Value = 7083.7896
$ws.Range("I136").Value = 7366.222
$ws.Range("K136").Value = 22098.666
$ws.Range("M136").Value = -19548.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4778.9443
$ws.Range("I134").Value = 5740.28
$ws.Range("J134").Value = 2594.0908
$ws.Range("K134").Value = 17220.84
$ws.Range("L134").Value = 7782.2724
$ws.Range("M134").Value = -14685.84
$ws.Range("N134").Value = -12852.2724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4442.7915
$ws.Range("I31").Value = 1630.6451
$ws.Range("J31").Value = 9570.823
$ws.Range("K31").Value = 1630.6451
$ws.Range("L31").Value = 9570.823
$ws.Range("M31").Value = -1335.6451
$ws.Range("N31").Value = -10160.823
$ws.Range("H34").Value = 4442.7915
$ws.Range("I34").Value = 1630.6451
$ws.Range("J34").Value = 9570.823
$ws.Range("K34").Value = 1630.6451
$ws.Range("L34").Value = 9570.823
$ws.Range("M34").Value = -1428.6451
$ws.Range("N34").Value = -9974.823
$ws.Range("H58").Value = 1535.9445
$ws.Range("I58").Value = 1235.4546
$ws.Range("J58").Value = 2008.1428
$ws.Range("K58").Value = 1235.4546
$ws.Range("L58").Value = 2008.1428
$ws.Range("M58").Value = -1032.4546
$ws.Range("N58").Value = -2414.1428
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 2918.6897
$ws.Range("I132").Value = 2984.0908
$ws.Range("J132").Value = 2713.1428
$ws.Range("K132").Value = 8952.2724
$ws.Range("L132").Value = 8139.428400000001
$ws.Range("M132").Value = -6422.2724
$ws.Range("N132").Value = -13199.4284
$ws.Range("H134").Value = 5029.3076
$ws.Range("I134").Value = 6067.2
$ws.Range("J134").Value = 1569.6666
$ws.Range("K134").Value = 18201.6
$ws.Range("L134").Value = 4708.9998
$ws.Range("M134").Value = -15666.6
$ws.Range("N134").Value = -9778.9998
$ws.Range("H136").Value = 1535.9445
$ws.Range("I136").Value = 1235.4546
$ws.Range("J136").Value = 2008.1428
$ws.Range("K136").Value = 3706.3638
$ws.Range("L136").Value = 6024.428400000001
$ws.Range("M136").Value = -1156.3638
$ws.Range("N136").Value = -11124.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 4513
$ws.Range("I130").Value = 965
$ws.Range("J130").Value = 5400
$ws.Range("K130").Value = 2895
$ws.Range("L130").Value = 16200
$ws.Range("M130").Value = 2125
$ws.Range("N130").Value = -26240
$ws.Range("H131").Value = 1471510.9
$ws.Range("I131").Value = 5556178
$ws.Range("J131").Value = 1030.58
$ws.Range("K131").Value = 16668534
$ws.Range("L131").Value = 3091.74
$ws.Range("M131").Value = -16663494
$ws.Range("N131").Value = -13171.74

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 70009
$ws.Range("J25").Value = 70009
$ws.Range("L25").Value = 70009
$ws.Range("N25").Value = -71067
$ws.Range("H103").Value = 43650
$ws.Range("J103").Value = 43650
$ws.Range("L103").Value = 43650
$ws.Range("N103").Value = -45994
$ws.Range("H132").Value = 3911.4412
$ws.Range("I132").Value = 4958.6875
$ws.Range("J132").Value = 2980.5557
$ws.Range("K132").Value = 14876.0625
$ws.Range("L132").Value = 8941.667099999999
$ws.Range("M132").Value = -12346.0625
$ws.Range("N132").Value = -14001.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1659722
$ws.Range("I22").Value = 4444915
$ws.Range("J22").Value = 1868.881
$ws.Range("K22").Value = 4444915
$ws.Range("L22").Value = 1868.881
$ws.Range("M22").Value = -4444620
$ws.Range("N22").Value = -2458.881
$ws.Range("H27").Value = 1659722
$ws.Range("I27").Value = 4444915
$ws.Range("J27").Value = 1868.881
$ws.Range("K27").Value = 4444915
$ws.Range("L27").Value = 1868.881
$ws.Range("M27").Value = -4444808
$ws.Range("N27").Value = -2082.881
$ws.Range("H46").Value = 15152443
$ws.Range("I46").Value = 30303658
$ws.Range("J46").Value = 1227
$ws.Range("K46").Value = 30303658
$ws.Range("L46").Value = 1227
$ws.Range("M46").Value = -30303470
$ws.Range("N46").Value = -1603
$ws.Range("H55").Value = 115384920
$ws.Range("I55").Value = 125000290
$ws.Range("J55").Value = 100000330
$ws.Range("K55").Value = 125000290
$ws.Range("L55").Value = 100000330
$ws.Range("M55").Value = -125000117
$ws.Range("N55").Value = -100000676
$ws.Range("H132").Value = 15283734
$ws.Range("I132").Value = 24130358
$ws.Range("J132").Value = 3200.7273
$ws.Range("K132").Value = 72391074
$ws.Range("L132").Value = 9602.1819
$ws.Range("M132").Value = -72388544
$ws.Range("N132").Value = -14662.1819
$ws.Range("H136").Value = 7352.067
$ws.Range("I136").Value = 6920.1
$ws.Range("J136").Value = 8216
$ws.Range("K136").Value = 20760.3
$ws.Range("L136").Value = 24648
$ws.Range("M136").Value = -18210.3
$ws.Range("N136").Value = -29748

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1583.081
$ws.Range("I132").Value = 1103
$ws.Range("J132").Value = 2287.2
$ws.Range("K132").Value = 3309
$ws.Range("L132").Value = 6861.599999999999
$ws.Range("M132").Value = -779
$ws.Range("N132").Value = -11921.6
$ws.Range("H136").Value = 3612.7727
$ws.Range("I136").Value = 4575.077
$ws.Range("J136").Value = 2222.7778
$ws.Range("K136").Value = 13725.231
$ws.Range("L136").Value = 6668.3334
$ws.Range("M136").Value = -11175.231
$ws.Range("N136").Value = -11768.3334
